$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("B4").Value = 168
$ws.Range("C4").Value = 196
$ws.Range("D4").Value = 134
$ws.Range("E4").Value = 148
$ws.Range("F4").Value = 114
$ws.Range("G4").Value = 123
$ws.Range("H4").Value = 102
$ws.Range("I4").Value = 108
$ws.Range("B5").Value = 158
$ws.Range("C5").Value = 205
$ws.Range("D5").Value = 129
$ws.Range("E5").Value = 154
$ws.Range("F5").Value = 112
$ws.Range("G5").Value = 128
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 111
$ws.Range("B6").Value = 126
$ws.Range("C6").Value = 243
$ws.Range("D6").Value = 112
$ws.Range("E6").Value = 181
$ws.Range("F6").Value = 101
$ws.Range("G6").Value = 147
$ws.Range("H6").Value = 94
$ws.Range("I6").Value = 126
$ws.Range("B7").Value = 127
$ws.Range("C7").Value = 224
$ws.Range("D7").Value = 112
$ws.Range("E7").Value = 168
$ws.Range("F7").Value = 101
$ws.Range("G7").Value = 138
$ws.Range("H7").Value = 93
$ws.Range("I7").Value = 119
$ws.Range("B8").Value = 140
$ws.Range("C8").Value = 228
$ws.Range("D8").Value = 120
$ws.Range("E8").Value = 170
$ws.Range("F8").Value = 106
$ws.Range("G8").Value = 138
$ws.Range("H8").Value = 97
$ws.Range("I8").Value = 119
$ws.Range("B9").Value = 181
$ws.Range("C9").Value = 204
$ws.Range("D9").Value = 142
$ws.Range("E9").Value = 154
$ws.Range("F9").Value = 119
$ws.Range("G9").Value = 127
$ws.Range("H9").Value = 105
$ws.Range("I9").Value = 110

$ws = $wb.Worksheets.Item(2)
$ws.Range("B4").Value = 300436
$ws.Range("C4").Value = 427218
$ws.Range("D4").Value = 314778
$ws.Range("E4").Value = 443031
$ws.Range("F4").Value = 315289
$ws.Range("G4").Value = 439106
$ws.Range("H4").Value = 311102
$ws.Range("I4").Value = 427606
$ws.Range("B5").Value = 307688
$ws.Range("C5").Value = 539275
$ws.Range("D5").Value = 316621
$ws.Range("E5").Value = 548079
$ws.Range("F5").Value = 314001
$ws.Range("G5").Value = 534708
$ws.Range("H5").Value = 307137
$ws.Range("I5").Value = 512932
$ws.Range("B6").Value = 251093
$ws.Range("C6").Value = 851383
$ws.Range("D6").Value = 241932
$ws.Range("E6").Value = 820018
$ws.Range("F6").Value = 232006
$ws.Range("G6").Value = 775169
$ws.Range("H6").Value = 221293
$ws.Range("I6").Value = 725033
$ws.Range("B7").Value = 110353
$ws.Range("C7").Value = 343477
$ws.Range("D7").Value = 107725
$ws.Range("E7").Value = 333150
$ws.Range("F7").Value = 103444
$ws.Range("G7").Value = 315626
$ws.Range("H7").Value = 98883
$ws.Range("I7").Value = 295683
$ws.Range("B8").Value = 47042
$ws.Range("C8").Value = 119606
$ws.Range("D8").Value = 46680
$ws.Range("E8").Value = 117731
$ws.Range("F8").Value = 45388
$ws.Range("G8").Value = 112546
$ws.Range("H8").Value = 43677
$ws.Range("I8").Value = 106273
$ws.Range("B9").Value = 20279439
$ws.Range("C9").Value = 26300578
$ws.Range("D9").Value = 21391712
$ws.Range("E9").Value = 27497039
$ws.Range("F9").Value = 21566370
$ws.Range("G9").Value = 27415896
$ws.Range("H9").Value = 21215148
$ws.Range("I9").Value = 26725005

$ws = $wb.Worksheets.Item(3)
$ws.Range("B4").Value = 113527
$ws.Range("C4").Value = 162645
$ws.Range("D4").Value = 114376
$ws.Range("E4").Value = 162047
$ws.Range("F4").Value = 110529
$ws.Range("G4").Value = 154537
$ws.Range("H4").Value = 104431
$ws.Range("I4").Value = 144045
$ws.Range("B5").Value = 114071
$ws.Range("C5").Value = 202439
$ws.Range("D5").Value = 112248
$ws.Range("E5").Value = 196833
$ws.Range("F5").Value = 106862
$ws.Range("G5").Value = 184179
$ws.Range("H5").Value = 100500
$ws.Range("I5").Value = 169607
$ws.Range("B6").Value = 85795
$ws.Range("C6").Value = 304111
$ws.Range("D6").Value = 79888
$ws.Range("E6").Value = 281849
$ws.Range("F6").Value = 73790
$ws.Range("G6").Value = 256056
$ws.Range("H6").Value = 68593
$ws.Range("I6").Value = 231073
$ws.Range("B7").Value = 38741
$ws.Range("C7").Value = 125305
$ws.Range("D7").Value = 36323
$ws.Range("E7").Value = 116676
$ws.Range("F7").Value = 33565
$ws.Range("G7").Value = 106207
$ws.Range("H7").Value = 31038
$ws.Range("I7").Value = 95961
$ws.Range("B8").Value = 16982
$ws.Range("C8").Value = 44334
$ws.Range("D8").Value = 16174
$ws.Range("E8").Value = 41829
$ws.Range("F8").Value = 15189
$ws.Range("G8").Value = 38461
$ws.Range("H8").Value = 14097
$ws.Range("I8").Value = 34964
$ws.Range("B9").Value = 7126521
$ws.Range("C9").Value = 9278465
$ws.Range("D9").Value = 7239374
$ws.Range("E9").Value = 9340646
$ws.Range("F9").Value = 7031779
$ws.Range("G9").Value = 8972138
$ws.Range("H9").Value = 6649480
$ws.Range("I9").Value = 8397334

$ws = $wb.Worksheets.Item(4)
$ws.Range("B4").Value = 6107
$ws.Range("C4").Value = 8771
$ws.Range("D4").Value = 6189
$ws.Range("E4").Value = 8790
$ws.Range("F4").Value = 6015
$ws.Range("G4").Value = 8431
$ws.Range("H4").Value = 5731
$ws.Range("I4").Value = 7939
$ws.Range("B5").Value = 5868
$ws.Range("C5").Value = 10438
$ws.Range("D5").Value = 5789
$ws.Range("E5").Value = 10172
$ws.Range("F5").Value = 5542
$ws.Range("G5").Value = 9557
$ws.Range("H5").Value = 5203
$ws.Range("I5").Value = 8824
$ws.Range("B6").Value = 4646
$ws.Range("C6").Value = 16368
$ws.Range("D6").Value = 4386
$ws.Range("E6").Value = 15305
$ws.Range("F6").Value = 4112
$ws.Range("G6").Value = 14048
$ws.Range("H6").Value = 3815
$ws.Range("I6").Value = 12790
$ws.Range("B7").Value = 2129
$ws.Range("C7").Value = 6852
$ws.Range("D7").Value = 2013
$ws.Range("E7").Value = 6439
$ws.Range("F7").Value = 1892
$ws.Range("G7").Value = 5934
$ws.Range("H7").Value = 1774
$ws.Range("I7").Value = 5431
$ws.Range("B8").Value = 924
$ws.Range("C8").Value = 2413
$ws.Range("D8").Value = 889
$ws.Range("E8").Value = 2298
$ws.Range("F8").Value = 839
$ws.Range("G8").Value = 2132
$ws.Range("H8").Value = 789
$ws.Range("I8").Value = 1962
$ws.Range("B9").Value = 376072
$ws.Range("C9").Value = 490583
$ws.Range("D9").Value = 384129
$ws.Range("E9").Value = 496497
$ws.Range("F9").Value = 374545
$ws.Range("G9").Value = 479648
$ws.Range("H9").Value = 358667
$ws.Range("I9").Value = 453962

$ws = $wb.Worksheets.Item(5)
$ws.Range("B4").Value = 2667
$ws.Range("C4").Value = 3818
$ws.Range("D4").Value = 2684
$ws.Range("E4").Value = 3805
$ws.Range("F4").Value = 2589
$ws.Range("G4").Value = 3620
$ws.Range("H4").Value = 2443
$ws.Range("I4").Value = 3375
$ws.Range("B5").Value = 2681
$ws.Range("C5").Value = 4755
$ws.Range("D5").Value = 2639
$ws.Range("E5").Value = 4623
$ws.Range("F5").Value = 2514
$ws.Range("G5").Value = 4326
$ws.Range("H5").Value = 2350
$ws.Range("I5").Value = 3975
$ws.Range("B6").Value = 2011
$ws.Range("C6").Value = 7133
$ws.Range("D6").Value = 1885
$ws.Range("E6").Value = 6617
$ws.Range("F6").Value = 1746
$ws.Range("G6").Value = 6014
$ws.Range("H6").Value = 1602
$ws.Range("I6").Value = 5418
$ws.Range("B7").Value = 906
$ws.Range("C7").Value = 2937
$ws.Range("D7").Value = 851
$ws.Range("E7").Value = 2735
$ws.Range("F7").Value = 793
$ws.Range("G7").Value = 2494
$ws.Range("H7").Value = 733
$ws.Range("I7").Value = 2251
$ws.Range("B8").Value = 398
$ws.Range("C8").Value = 1040
$ws.Range("D8").Value = 380
$ws.Range("E8").Value = 982
$ws.Range("F8").Value = 355
$ws.Range("G8").Value = 902
$ws.Range("H8").Value = 330
$ws.Range("I8").Value = 819
$ws.Range("B9").Value = 167136
$ws.Range("C9").Value = 217694
$ws.Range("D9").Value = 170037
$ws.Range("E9").Value = 219396
$ws.Range("F9").Value = 164709
$ws.Range("G9").Value = 210201
$ws.Range("H9").Value = 156344
$ws.Range("I9").Value = 197388

$ws = $wb.Worksheets.Item(6)
$ws.Range("B4").Value = 2087480
$ws.Range("C4").Value = 3184944
$ws.Range("D4").Value = 1655416
$ws.Range("E4").Value = 2389489
$ws.Range("F4").Value = 1337758
$ws.Range("G4").Value = 1878464
$ws.Range("H4").Value = 1095395
$ws.Range("I4").Value = 1503909
$ws.Range("B5").Value = 1907856
$ws.Range("C5").Value = 3782354
$ws.Range("D5").Value = 1516774
$ws.Range("E5").Value = 2767133
$ws.Range("F5").Value = 1223920
$ws.Range("G5").Value = 2143974
$ws.Range("H5").Value = 999508
$ws.Range("I5").Value = 1692830
$ws.Range("B6").Value = 991914
$ws.Range("C6").Value = 4824685
$ws.Range("D6").Value = 807889
$ws.Range("E6").Value = 3334270
$ws.Range("F6").Value = 666287
$ws.Range("G6").Value = 2565836
$ws.Range("H6").Value = 555775
$ws.Range("I6").Value = 2016151
$ws.Range("B7").Value = 466365
$ws.Range("C7").Value = 1894532
$ws.Range("D7").Value = 378593
$ws.Range("E7").Value = 1394784
$ws.Range("F7").Value = 311423
$ws.Range("G7").Value = 1074002
$ws.Range("H7").Value = 259241
$ws.Range("I7").Value = 844640
$ws.Range("B8").Value = 231340
$ws.Range("C8").Value = 759499
$ws.Range("D8").Value = 186498
$ws.Range("E8").Value = 531223
$ws.Range("F8").Value = 152441
$ws.Range("G8").Value = 409563
$ws.Range("H8").Value = 126139
$ws.Range("I8").Value = 322523
$ws.Range("B9").Value = 154002774
$ws.Range("C9").Value = 238957293
$ws.Range("D9").Value = 114786894
$ws.Range("E9").Value = 149958184
$ws.Range("F9").Value = 92154516
$ws.Range("G9").Value = 117815280
$ws.Range("H9").Value = 75006315
$ws.Range("I9").Value = 94264781

$ws = $wb.Worksheets.Item(7)
$ws.Range("B4").Value = 323594
$ws.Range("C4").Value = 493013
$ws.Range("D4").Value = 268217
$ws.Range("E4").Value = 384473
$ws.Range("F4").Value = 225714
$ws.Range("G4").Value = 315319
$ws.Range("H4").Value = 191838
$ws.Range("I4").Value = 262426
$ws.Range("B5").Value = 298011
$ws.Range("C5").Value = 587387
$ws.Range("D5").Value = 246553
$ws.Range("E5").Value = 443266
$ws.Range("F5").Value = 206151
$ws.Range("G5").Value = 356777
$ws.Range("H5").Value = 173934
$ws.Range("I5").Value = 291631
$ws.Range("B6").Value = 169980
$ws.Range("C6").Value = 787655
$ws.Range("D6").Value = 143766
$ws.Range("E6").Value = 557425
$ws.Range("F6").Value = 122881
$ws.Range("G6").Value = 447340
$ws.Range("H6").Value = 106094
$ws.Range("I6").Value = 365624
$ws.Range("B7").Value = 79529
$ws.Range("C7").Value = 306010
$ws.Range("D7").Value = 67127
$ws.Range("E7").Value = 234568
$ws.Range("F7").Value = 57300
$ws.Range("G7").Value = 188515
$ws.Range("H7").Value = 49442
$ws.Range("I7").Value = 154369
$ws.Range("B8").Value = 38188
$ws.Range("C8").Value = 123009
$ws.Range("D8").Value = 32058
$ws.Range("E8").Value = 88153
$ws.Range("F8").Value = 27210
$ws.Range("G8").Value = 70913
$ws.Range("H8").Value = 23342
$ws.Range("I8").Value = 58117
$ws.Range("B9").Value = 22704226
$ws.Range("C9").Value = 36281378
$ws.Range("D9").Value = 17341947
$ws.Range("E9").Value = 22542049
$ws.Range("F9").Value = 14557796
$ws.Range("G9").Value = 18536838
$ws.Range("H9").Value = 12341708
$ws.Range("I9").Value = 15465169

$ws = $wb.Worksheets.Item(8)
$ws.Range("B4").Value = 16078
$ws.Range("C4").Value = 24505
$ws.Range("D4").Value = 13808
$ws.Range("E4").Value = 19843
$ws.Range("F4").Value = 12066
$ws.Range("G4").Value = 16930
$ws.Range("H4").Value = 10672
$ws.Range("I4").Value = 14686
$ws.Range("B5").Value = 14215
$ws.Range("C5").Value = 27893
$ws.Range("D5").Value = 12117
$ws.Range("E5").Value = 21739
$ws.Range("F5").Value = 10460
$ws.Range("G5").Value = 18095
$ws.Range("H5").Value = 9131
$ws.Range("I5").Value = 15328
$ws.Range("B6").Value = 9181
$ws.Range("C6").Value = 40363
$ws.Range("D6").Value = 8094
$ws.Range("E6").Value = 30003
$ws.Range("F6").Value = 7221
$ws.Range("G6").Value = 25178
$ws.Range("H6").Value = 6514
$ws.Range("I6").Value = 21549
$ws.Range("B7").Value = 4332
$ws.Range("C7").Value = 16012
$ws.Range("D7").Value = 3819
$ws.Range("E7").Value = 12868
$ws.Range("F7").Value = 3411
$ws.Range("G7").Value = 10846
$ws.Range("H7").Value = 3082
$ws.Range("I7").Value = 9328
$ws.Range("B8").Value = 2009
$ws.Range("C8").Value = 6316
$ws.Range("D8").Value = 1757
$ws.Range("E8").Value = 4746
$ws.Range("F8").Value = 1556
$ws.Range("G8").Value = 3996
$ws.Range("H8").Value = 1396
$ws.Range("I8").Value = 3433
$ws.Range("B9").Value = 1102993
$ws.Range("C9").Value = 1759791
$ws.Range("D9").Value = 870700
$ws.Range("E9").Value = 1134234
$ws.Range("F9").Value = 757691
$ws.Range("G9").Value = 968127
$ws.Range("H9").Value = 667214
$ws.Range("I9").Value = 839943

$ws = $wb.Worksheets.Item(9)
$ws.Range("B4").Value = 76755
$ws.Range("C4").Value = 120319
$ws.Range("D4").Value = 63046
$ws.Range("E4").Value = 90400
$ws.Range("F4").Value = 53049
$ws.Range("G4").Value = 74110
$ws.Range("H4").Value = 45087
$ws.Range("I4").Value = 61677
$ws.Range("B5").Value = 70280
$ws.Range("C5").Value = 144921
$ws.Range("D5").Value = 57950
$ws.Range("E5").Value = 104240
$ws.Range("F5").Value = 48451
$ws.Range("G5").Value = 83854
$ws.Range("H5").Value = 40879
$ws.Range("I5").Value = 68542
$ws.Range("B6").Value = 39951
$ws.Range("C6").Value = 200504
$ws.Range("D6").Value = 33789
$ws.Range("E6").Value = 131181
$ws.Range("F6").Value = 28881
$ws.Range("G6").Value = 105140
$ws.Range("H6").Value = 24935
$ws.Range("I6").Value = 85932
$ws.Range("B7").Value = 18692
$ws.Range("C7").Value = 73972
$ws.Range("D7").Value = 15777
$ws.Range("E7").Value = 55148
$ws.Range("F7").Value = 13467
$ws.Range("G7").Value = 44306
$ws.Range("H7").Value = 11620
$ws.Range("I7").Value = 36281
$ws.Range("B8").Value = 8978
$ws.Range("C8").Value = 31282
$ws.Range("D8").Value = 7535
$ws.Range("E8").Value = 20742
$ws.Range("F8").Value = 6395
$ws.Range("G8").Value = 16667
$ws.Range("H8").Value = 5486
$ws.Range("I8").Value = 13659
$ws.Range("B9").Value = 5891822
$ws.Range("C9").Value = 10087279
$ws.Range("D9").Value = 4081402
$ws.Range("E9").Value = 5316173
$ws.Range("F9").Value = 3421600
$ws.Range("G9").Value = 4356959
$ws.Range("H9").Value = 2900653
$ws.Range("I9").Value = 3634759
